$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145, shifting existing rows 145:177 down to 146:178
$ws.Rows("145:145").Insert()

# Populate the newly inserted row 145 with the new weekly data record
$ws.Range("A145").Value = 6
$ws.Range("B145").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C145").Value = "Metropolitana"
$ws.Range("D145").Value = 44627
$ws.Range("E145").Value = 13
$ws.Range("F145").Value = "Fruta"
$ws.Range("G145").Value = 100101
$ws.Range("H145").Value = "Berries"
$ws.Range("I145").Value = 100101004
$ws.Range("J145").Value = "Frambuesa"
$ws.Range("K145").Value = "Sin especificar"
$ws.Range("L145").Value = "Especial"
$ws.Range("M145").Value = 200
$ws.Range("N145").Value = 8000
$ws.Range("O145").Value = 8000
$ws.Range("P145").Value = 8000
$ws.Range("Q145").Value = '$/bandeja 2 kilos'
$ws.Range("R145").Value = "Provincia de Linares"
$ws.Range("S145").Value = 4000
$ws.Range("T145").Value = 2
